$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("PLC Tags")

# Update the "Logical Address" column (D) for the PLC output tags whose
# addresses were renumbered from %Q1.x / %Q2.x to the new %Q2.x scheme.
$ws.Range("D15").Value = "%Q2.0"
$ws.Range("D16").Value = "%Q2.1"
$ws.Range("D17").Value = "%Q2.2"
$ws.Range("D18").Value = "%Q2.3"
$ws.Range("D19").Value = "%Q2.4"
